# Update the dSF column (F) with re-pulled / recalculated data.
# Rows 11 and 18 are unchanged; all other data rows (2-23) get new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 1
    3  = 1
    4  = -4
    5  = -4
    6  = -2
    7  = -5
    8  = -6
    9  = 4
    10 = -1
    12 = 3
    13 = 3
    14 = -1
    15 = -6
    16 = -3
    17 = -3
    19 = 1
    20 = -2
    21 = 1
    22 = 2
    23 = -6
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 6).Value = $values[$row]
}
